# edits for JWM revision.  changed gamma to phi to avoid confusion with "y".
# This was done in the text and sup info, but not the actual model code!

$wb = $excel.ActiveWorkbook

# Rename the "for R" sheet to "for stat analysis"
$wsForR = $wb.Worksheets.Item("for R")
$wsForR.Name = "for stat analysis"

# Fix the shared formula range in column N (N533:N596 -> N533:N576)
# and update the active selection on this sheet.
$wsForR.Range("N533:N576").FormulaR1C1 = "=RC13-RC7"
$wsForR.Range("F11").Select()

# Sheet "1-769": remove scrolled topLeftCell, restore view to top
$ws1769 = $wb.Worksheets.Item("1-769")
$ws1769.Activate()
$ws1769.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Sheet "770-825": remove scrolled topLeftCell, restore view to top
$ws770825 = $wb.Worksheets.Item("770-825")
$ws770825.Activate()
$ws770825.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Re-activate "for stat analysis" sheet, matching tabSelected state
$wsForR.Activate()
